$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F3").Value = -2.283908985393195
$ws.Range("F4").Value = -2.416088453551624
$ws.Range("F5").Value = -2.548967639398655
$ws.Range("F6").Value = -2.681846825245684
$ws.Range("F7").Value = -2.814726011092715
$ws.Range("F8").Value = -2.948304914628344
$ws.Range("F9").Value = -3.081184100475374
$ws.Range("F10").Value = -3.214763004011005
$ws.Range("F11").Value = -3.348341907546635
$ws.Range("F12").Value = -3.481920811082266
$ws.Range("F13").Value = -3.616199432306497
$ws.Range("F14").Value = -3.75117777121933
$ws.Range("F15").Value = -3.855339157786533
$ws.Range("F17").Value = -2.497660976620755
$ws.Range("F18").Value = -2.529112362948664
$ws.Range("F19").Value = -2.560563749276572
$ws.Range("F20").Value = -2.59201513560448
$ws.Range("F21").Value = -2.623466521932388
$ws.Range("F22").Value = -2.654917908260296
$ws.Range("F23").Value = -2.686369294588204
$ws.Range("F24").Value = -2.717820680916112
$ws.Range("F25").Value = -2.74927206724402
$ws.Range("F26").Value = -2.780723453571928
$ws.Range("F27").Value = -2.812174839899834
$ws.Range("F28").Value = -2.843626226227742
$ws.Range("F29").Value = -2.873173192687904
$ws.Range("F30").Value = -2.902720159148066
$ws.Range("F31").Value = -2.932267125608218
$ws.Range("F32").Value = -2.961814092068379
$ws.Range("F33").Value = -2.991361058528541
$ws.Range("F34").Value = -3.020908024988696
$ws.Range("F35").Value = -3.050454991448855
$ws.Range("F36").Value = -3.080001957909011
$ws.Range("F37").Value = -3.109056859205772
$ws.Range("F38").Value = -3.13761969533912
$ws.Range("F39").Value = -3.167166661799279
$ws.Range("F40").Value = -3.196221563096033
$ws.Range("F41").Value = -3.224292334065978
$ws.Range("F42").Value = -3.253347235362731
$ws.Range("F43").Value = -3.281910071496083
$ws.Range("F44").Value = -3.310472907629436
$ws.Range("F45").Value = -3.339035743762789
$ws.Range("F46").Value = -3.367106514732737
$ws.Range("F47").Value = -3.395669350866083
$ws.Range("F48").Value = -3.424232186999435
$ws.Range("F49").Value = -3.452302957969381
$ws.Range("F50").Value = -3.479389598612523
$ws.Range("F51").Value = -3.50746036958247
$ws.Range("F52").Value = -3.535039075389013
$ws.Range("F53").Value = -3.562617781195557
$ws.Range("F54").Value = -3.590196487002101
$ws.Range("F55").Value = -3.617775192808644
$ws.Range("F56").Value = -3.645353898615188
$ws.Range("F57").Value = -3.672932604421731
$ws.Range("F58").Value = -3.700511310228275
$ws.Range("F59").Value = -3.726649399536589
$ws.Range("F60").Value = -3.756203857030933
$ws.Range("F61").Value = -3.78170163320847
$ws.Range("F62").Value = -3.806559096255227
$ws.Range("F63").Value = -3.831416559301985
$ws.Range("F64").Value = -3.856274022348742
$ws.Range("F65").Value = -3.880491172264722
$ws.Range("F66").Value = -3.904068009049922
$ws.Range("F67").Value = -3.927644845835124
$ws.Range("F68").Value = -3.950581369489547
$ws.Range("F69").Value = -3.972877580013191
$ws.Range("F70").Value = -3.995173790536836
$ws.Range("F71").Value = -4.017470001060481
$ws.Range("F309").Value = -3.018311655611262
$ws.Range("F310").Value = -3.17845730755005
$ws.Range("F311").Value = -3.334234127884221
$ws.Range("F312").Value = -3.486326736666184
$ws.Range("F313").Value = -3.634333330164532
$ws.Range("F314").Value = -3.780268381506738
$ws.Range("F315").Value = -3.922024566550103
$ws.Range("F316").Value = -4.057763571624702
$ws.Range("F317").Value = -4.188725883400633
$ws.Range("F318").Value = -4.314433232740556

Write-Output "Updated 78 cells"
